$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

# Rows 12-14 (metrics enterpriseCount, birthRate, deathRate) share the same
# dataText string describing what an "enterprise" is. Reword it to talk
# about "businesses" instead, per the commit message "Change name to
# businesses".
$newText = "Represented here are enterprises, which can be thought of as the overall business, made up of all the individual sites or workplaces."

$ws.Range("C12").Value = $newText
$ws.Range("C13").Value = $newText
$ws.Range("C14").Value = $newText

# Reflect the updated selection/view state from the saved workbook.
$ws.Activate() | Out-Null
$ws.Range("C12:C14").Select() | Out-Null

$wb.Save()
